$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row total correct answers: 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row correct marks: 69 -> 115
$ws.Range("B12").Value = 115

# Update the Corr/total marks text: 66/84 -> 115/140
$ws.Range("E12").Value = "115/140"
